$wb = $excel.ActiveWorkbook

$women = $wb.Worksheets.Item("Women")
$men = $wb.Worksheets.Item("Men")

# ---- Sheet "Men": rows 2-4 first, to control shared-string insertion order ----

# Row 2
$men.Cells.Item(2, 1).Value = 5
$men.Cells.Item(2, 2).Value = "hehe"
$men.Cells.Item(2, 3).Value = "fefe"
$men.Cells.Item(2, 4).Value = 0
$men.Cells.Item(2, 5).Value = 0
$men.Cells.Item(2, 6).Value = 0
$men.Cells.Item(2, 7).Value = 0
$men.Cells.Item(2, 8).Value = 0
$men.Cells.Item(2, 9).Value = 0
$men.Cells.Item(2, 10).Value = 0
$men.Cells.Item(2, 11).Value = 0
$men.Cells.Item(2, 12).Value = 0
$men.Cells.Item(2, 13).Value = 0

# Row 3
$men.Cells.Item(3, 1).Value = 6
$men.Cells.Item(3, 2).Value = "john"
$men.Cells.Item(3, 3).Value = "kokokaka"
$men.Cells.Item(3, 4).Value = 0
$men.Cells.Item(3, 5).Value = 0
$men.Cells.Item(3, 6).Value = 0
$men.Cells.Item(3, 7).Value = 0
$men.Cells.Item(3, 8).Value = 0
$men.Cells.Item(3, 9).Value = 0
$men.Cells.Item(3, 10).Value = 0
$men.Cells.Item(3, 11).Value = 0
$men.Cells.Item(3, 12).Value = 0
$men.Cells.Item(3, 13).Value = 0

# Row 4
$men.Cells.Item(4, 1).Value = 7
$men.Cells.Item(4, 2).Value = "jahny"
$men.Cells.Item(4, 3).Value = "hayayay"
$men.Cells.Item(4, 4).Value = 929
$men.Cells.Item(4, 5).Value = 461
$men.Cells.Item(4, 6).Value = 1181
$men.Cells.Item(4, 7).Value = 901
$men.Cells.Item(4, 8).Value = 751
$men.Cells.Item(4, 9).Value = 825
$men.Cells.Item(4, 10).Value = 1321
$men.Cells.Item(4, 11).Value = 509
$men.Cells.Item(4, 12).Value = 972
$men.Cells.Item(4, 13).Value = 264

# ---- Sheet "Women": add row 6 ----

$women.Cells.Item(6, 1).Value = 8
$women.Cells.Item(6, 2).Value = "woman"
$women.Cells.Item(6, 3).Value = "womansson"
$women.Cells.Item(6, 4).Value = 690
$women.Cells.Item(6, 5).Value = 1420
$women.Cells.Item(6, 6).Value = 995
$women.Cells.Item(6, 7).Value = 584
$women.Cells.Item(6, 8).Value = 1082
$women.Cells.Item(6, 9).Value = 1094
$women.Cells.Item(6, 10).Value = 1233

# ---- Sheet "Men": rows 5-6 ----

# Row 5
$men.Cells.Item(5, 1).Value = 9
$men.Cells.Item(5, 2).Value = "john"
$men.Cells.Item(5, 3).Value = "smithy"
$men.Cells.Item(5, 4).Value = 929
$men.Cells.Item(5, 5).Value = 461
$men.Cells.Item(5, 6).Value = 1181
$men.Cells.Item(5, 7).Value = 901
$men.Cells.Item(5, 8).Value = 751
$men.Cells.Item(5, 9).Value = 825
$men.Cells.Item(5, 10).Value = 771
$men.Cells.Item(5, 11).Value = 509
$men.Cells.Item(5, 12).Value = 972
$men.Cells.Item(5, 13).Value = 264

# Row 6
$men.Cells.Item(6, 1).Value = 10
$men.Cells.Item(6, 2).Value = "barack"
$men.Cells.Item(6, 3).Value = "obama"
$men.Cells.Item(6, 4).Value = 929
$men.Cells.Item(6, 5).Value = 461
$men.Cells.Item(6, 6).Value = 1181
$men.Cells.Item(6, 7).Value = 901
$men.Cells.Item(6, 8).Value = 751
$men.Cells.Item(6, 9).Value = 825
$men.Cells.Item(6, 10).Value = 895
$men.Cells.Item(6, 11).Value = 509
$men.Cells.Item(6, 12).Value = 972
$men.Cells.Item(6, 13).Value = 264

Write-Output "edit complete"
